$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.22"
$ws.Range("D3").Value = "'22.41"
$ws.Range("D4").Value = "'5.470"
$ws.Range("D7").Value = "'0.8041"
$ws.Range("D8").Value = "'1.049"
$ws.Range("D10").Value = "'0.07283"
$ws.Range("D11").Value = "'0.03178"
$ws.Range("D12").Value = "'0.02947"
$ws.Range("D13").Value = "'0.09257"
$ws.Range("D14").Value = "'0.001669"
$ws.Range("D15").Value = "'3.194"
$ws.Range("D16").Value = "'0.04697"
$ws.Range("D17").Value = "'0.01162"
$ws.Range("E17").Value = "16OneONEBestin24h"
$ws.Range("D18").Value = "'0.006279"
$ws.Range("D19").Value = "'0.001057"
$ws.Range("D20").Value = "'0.003815"
$ws.Range("D21").Value = "'0.0001504"
$ws.Range("D22").Value = "'0.0003608"
$ws.Range("D23").Value = "'3.977"
$ws.Range("D24").Value = "'3.395"
$ws.Range("D25").Value = "'2.118"
$ws.Range("E27").Value = "26ProBitTokenPROB"
$ws.Range("D40").Value = "'0.04155"
$ws.Range("D41").Value = "'0.006976"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003508"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1041"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("D44").Value = "'0.01026"
$ws.Range("D45").Value = "'0.00005642"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("D47").Value = "'0.6816"
$ws.Range("D48").Value = "'0.02604"
$ws.Range("D49").Value = "'0.00002105"
